$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.219.02"
$ws.Range("E2").Value = "  -2.94%  "

# Row 3
$ws.Range("D3").Value = "1.551.73"
$ws.Range("E3").Value = "  -4.73%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "'206.43"
$ws.Range("E5").Value = "  -3.55%  "

# Row 6
$ws.Range("E6").Value = "  +0.02%  "

# Row 7
$ws.Range("D7").Value = "'0.476"
$ws.Range("E7").Value = "  -5.34%  "

# Row 8
$ws.Range("D8").Value = "'0.0605"
$ws.Range("E8").Value = "  -2.08%  "

# Row 9
$ws.Range("D9").Value = "'0.241"
$ws.Range("E9").Value = "  -3.87%  "

# Row 10
$ws.Range("D10").Value = "'17.79"
$ws.Range("E10").Value = "  -3.26%  "

# Row 11
$ws.Range("D11").Value = "'0.0779"
$ws.Range("E11").Value = "  -1.23%  "

# Row 12
$ws.Range("D12").Value = "1.764.99"
$ws.Range("E12").Value = "  -4.76%  "

# Row 13
$ws.Range("D13").Value = "1.539.04"
$ws.Range("E13").Value = "  -6.92%  "

# Row 14
$ws.Range("D14").Value = "'3.98"
$ws.Range("E14").Value = "  -4.97%  "

# Row 15
$ws.Range("D15").Value = "'0.503"
$ws.Range("E15").Value = "  -4.61%  "

# Row 16
$ws.Range("D16").Value = "25.177.90"
$ws.Range("E16").Value = "  -2.97%  "

# Row 17
$ws.Range("D17").Value = "0.0₃0708"
$ws.Range("E17").Value = "  -4.11%  "

# Row 18
$ws.Range("D18").Value = "'58.59"
$ws.Range("E18").Value = "  -4.42%  "

# Row 19
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("D20").Value = "'184.86"

# Row 21
$ws.Range("D21").Value = "'4.10"
$ws.Range("E21").Value = "  -3.34%  "

# Row 22
$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -3.78%  "

# Row 23
$ws.Range("D23").Value = "'5.83"
$ws.Range("E23").Value = "  -4.01%  "

# Row 24
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("E25").Value = "  -3.87%  "

# Row 26
$ws.Range("D26").Value = "'139.26"
$ws.Range("E26").Value = "  -3.10%  "

# Row 27
$ws.Range("E27").Value = "  -4.58%  "

# Row 28
$ws.Range("D28").Value = "'14.76"
$ws.Range("E28").Value = "  -2.97%  "

# Row 29
$ws.Range("D29").Value = "'6.38"
$ws.Range("E29").Value = "  -5.15%  "

# Row 30
$ws.Range("D30").Value = "'1.15"
$ws.Range("E30").Value = "  -6.98%  "

# Row 31
$ws.Range("D31").Value = "'0.0461"
$ws.Range("E31").Value = "  -4.58%  "

# Row 32
$ws.Range("E32").Value = "  -3.70%  "

# Row 33
$ws.Range("E33").Value = "  -4.75%  "

# Row 34
$ws.Range("D34").Value = "'1.45"
$ws.Range("E34").Value = "  -3.83%  "

# Row 35
$ws.Range("E35").Value = "  -4.06%  "

# Row 36
$ws.Range("D36").Value = "1.086.60"
$ws.Range("E36").Value = "  -3.51%  "

# Row 37
$ws.Range("E37").Value = "  -0.41%  "

# Row 38
$ws.Range("D38").Value = "'0.0149"
$ws.Range("E38").Value = "  -2.67%  "

# Row 39
$ws.Range("D39").Value = "'0.492"
$ws.Range("E39").Value = "  -5.39%  "

# Row 40
$ws.Range("E40").Value = "  -7.49%  "

# Row 41
$ws.Range("D41").Value = "'0.759"
$ws.Range("E41").Value = "  -10.71%  "

# Row 42
$ws.Range("D42").Value = "'0.799"
$ws.Range("E42").Value = "  +4.73%  "

# Row 43
$ws.Range("D43").Value = "'92.42"
$ws.Range("E43").Value = "  -5.85%  "

# Row 44
$ws.Range("D44").Value = "'5.02"
$ws.Range("E44").Value = "  -2.76%  "

# Row 45
$ws.Range("D45").Value = "1.681.74"
$ws.Range("E45").Value = "  -4.62%  "

# Row 46
$ws.Range("D46").Value = "0.0₆0111"
$ws.Range("E46").Value = "  +5.73%  "

# Row 47
$ws.Range("D47").Value = "'52.23"
$ws.Range("E47").Value = "  -4.00%  "

# Row 48
$ws.Range("E48").Value = "  -2.33%  "

# Row 49
$ws.Range("E49").Value = "  -5.71%  "

# Row 50
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.23%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.404"
$ws.Range("E51").Value = "  -2.09%  "
